# Running test specifications, adding comments, and defining function for
# inflation target estimate.
#
# Adds two indicator columns to "sheet1" that flag which RBA Governor was
# in office for each monthly observation:
#   BU = "Glen Stevens"  (Governor until ~Sep 2016)
#   BV = "Phillip Lowe"  (Governor from ~Oct 2016)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sheet1")

$BU = 73   # column BU
$BV = 74   # column BV

# Header row (bold, matching the other header cells).
$ws.Cells.Item(1, $BU).Value = "Glen Stevens"
$ws.Cells.Item(1, $BU).Font.Bold = $true
$ws.Cells.Item(1, $BV).Value = "Phillip Lowe"
$ws.Cells.Item(1, $BV).Font.Bold = $true

# Data rows 2..161 -> Glen Stevens governed through the row at A112
# (Oct-2016); Phillip Lowe from then on.
$lastStevensRow = 111

for ($r = 2; $r -le 161; $r++) {
    if ($r -le $lastStevensRow) {
        $ws.Cells.Item($r, $BU).Value = 1
        $ws.Cells.Item($r, $BV).Value = 0
    } else {
        $ws.Cells.Item($r, $BU).Value = 0
        $ws.Cells.Item($r, $BV).Value = 1
    }
}

# Restore the view/selection as closely as possible: the workbook was left
# scrolled over to the new columns, with the active cell at BX15.
$ws.Range("F1").Select()
$ws.Range("BX15").Select()
